$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# The existing "Comment" column (AJ) moves to AM; three new GLP-1 columns
# are inserted at AJ:AL. Copy AJ1's header formatting (bold white text on
# navy fill, centered/wrapped) onto the newly used header cells first, so
# the look matches the rest of the header row, then fill in the text.
$ws.Range("AJ1").Copy()
$ws.Range("AK1:AM1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AM1").Value = "Comment"
$ws.Range("AJ1").Value = "GLP-1 Rx Rate"
$ws.Range("AK1").Value = "GLP-1 Market Penetration"
$ws.Range("AL1").Value = "GLP-1 Est Users"

# --- Data rows (2-22) ---------------------------------------------------
# For every hospital row: shift the existing "Comment" text out to the new
# AM column, then populate the three new GLP-1 metrics in AJ:AL.
$ws.Range("AM2").Value = "Regional CSC hub — only Level I trauma in SE Georgia"
$ws.Range("AJ2").Value = 58.4
$ws.Range("AK2").Value = 22.1
$ws.Range("AL2").Value = 14960
$ws.Range("AM3").Value = "AdventHealth system CSC campus — DNV certified 2019"
$ws.Range("AJ3").Value = 72.6
$ws.Range("AK3").Value = 26.8
$ws.Range("AL3").Value = 22780
$ws.Range("AM4").Value = "Flagship CSC — largest stroke program in NE Florida"
$ws.Range("AJ4").Value = 64.2
$ws.Range("AK4").Value = 24.3
$ws.Range("AL4").Value = 42350
$ws.Range("AM5").Value = "Academic flagship — 1,043 stroke cases in 2024 (AHA data)"
$ws.Range("AJ5").Value = 55.8
$ws.Range("AK5").Value = 21.4
$ws.Range("AL5").Value = 18620
$ws.Range("AM6").Value = "Academic med center — TraumaOne helicopter stroke transport"
$ws.Range("AJ6").Value = 62.1
$ws.Range("AK6").Value = 23.7
$ws.Range("AL6").Value = 33880
$ws.Range("AM7").Value = "Pepin Heart & Vascular Institute — BayCare system"
$ws.Range("AJ7").Value = 68.7
$ws.Range("AK7").Value = 25.6
$ws.Range("AL7").Value = 38250
$ws.Range("AM8").Value = "Ranked #3 in Jacksonville by U.S. News 2025"
$ws.Range("AJ8").Value = 61.8
$ws.Range("AK8").Value = 23.1
$ws.Range("AL8").Value = 19640
$ws.Range("AM9").Value = "HCA facility — CSC certified, serves 14 surrounding counties"
$ws.Range("AJ9").Value = 54.3
$ws.Range("AK9").Value = 20.8
$ws.Range("AL9").Value = 10080
$ws.Range("AM10").Value = "Only CSC between Jacksonville and Pensacola — serves 21 counties"
$ws.Range("AJ10").Value = 49.6
$ws.Range("AK10").Value = 19.2
$ws.Range("AL10").Value = 10560
$ws.Range("AM11").Value = "Satellite of UF Health — spoke to main Shands CSC"
$ws.Range("AJ11").Value = 52.1
$ws.Range("AK11").Value = 20.3
$ws.Range("AL11").Value = 9130
$ws.Range("AM12").Value = "Award-winning cardiac and stroke programs — CSC certified"
$ws.Range("AJ12").Value = 66.4
$ws.Range("AK12").Value = 25.1
$ws.Range("AL12").Value = 8820
$ws.Range("AM13").Value = "Treasure Coast — only Advanced Thrombectomy Center in St. Lucie County"
$ws.Range("AJ13").Value = 74.8
$ws.Range("AK13").Value = 27.9
$ws.Range("AL13").Value = 10920
$ws.Range("AM14").Value = "Panhandle Level II Trauma — nearest CSC is TMH (100+ mi)"
$ws.Range("AJ14").Value = 46.2
$ws.Range("AK14").Value = 18.4
$ws.Range("AL14").Value = 5150
$ws.Range("AM15").Value = "Ranked #1 in FL by U.S. News 2025 — expanded to 419 beds (2024)"
$ws.Range("AJ15").Value = 71.3
$ws.Range("AK15").Value = 26.4
$ws.Range("AL15").Value = 22140
$ws.Range("AM16").Value = "Opened 2017 — 7 neuro ORs incl. 2 intraoperative MRI suites"
$ws.Range("AJ16").Value = 53.7
$ws.Range("AK16").Value = 20.9
$ws.Range("AL16").Value = 12540
$ws.Range("AM17").Value = "Second-oldest US hospital in continuous operation (est. 1804)"
$ws.Range("AJ17").Value = 56.7
$ws.Range("AK17").Value = 21.8
$ws.Range("AL17").Value = 8550
$ws.Range("AM18").Value = "VA 1a High Complexity — very high 65+ population"
$ws.Range("AJ18").Value = 82.4
$ws.Range("AK18").Value = 31.6
$ws.Range("AL18").Value = 26040
$ws.Range("AM19").Value = "Thrombectomy-Capable cert renewed April 2025 — Level II Trauma"
$ws.Range("AJ19").Value = 76.3
$ws.Range("AK19").Value = 28.7
$ws.Range("AL19").Value = 17400
$ws.Range("AM20").Value = "America's 250 Best Hospitals Award (Healthgrades)"
$ws.Range("AJ20").Value = 60.5
$ws.Range("AK20").Value = 22.8
$ws.Range("AL20").Value = 16100
$ws.Range("AM21").Value = "Only Primary Stroke Center in Clay County — 6-county catchment"
$ws.Range("AJ21").Value = 57.3
$ws.Range("AK21").Value = 21.6
$ws.Range("AL21").Value = 5880
$ws.Range("AM22").Value = "3rd largest US hospital — `$660M expansion adding 440 beds (2025)"
$ws.Range("AJ22").Value = 69.8
$ws.Range("AK22").Value = 25.9
$ws.Range("AL22").Value = 54180
